$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text
# format first, so Excel stores them as text (matching the workbook's
# existing inline-string cells) instead of auto-converting to numeric values.
$textFormatCells = @("D5", "D6", "D8", "D10", "D12", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D30", "D32", "D33", "D35", "D36", "D38", "D40", "D41", "D42", "D44", "D46", "D48", "D50", "D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin values
$ws.Range("D2").Value = "68.708.10"
$ws.Range("E2").Value = "  +2.78%  "
$ws.Range("D3").Value = "2.653.71"
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "601.11"
$ws.Range("E5").Value = "  +2.52%  "
$ws.Range("D6").Value = "155.84"
$ws.Range("E6").Value = "  +5.15%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "2.649.40"
$ws.Range("E9").Value = "  +2.50%  "
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  +14.66%  "
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "5.24"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("E13").Value = "  +3.25%  "
$ws.Range("D14").Value = "27.98"
$ws.Range("E14").Value = "  +3.76%  "
$ws.Range("E15").Value = "  +7.41%  "
$ws.Range("D16").Value = "3.131.88"
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").Value = "68.542.34"
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("D18").Value = "2.654.54"
$ws.Range("E18").Value = "  +2.44%  "
$ws.Range("D19").Value = "11.39"
$ws.Range("E19").Value = "  +4.63%  "
$ws.Range("D20").Value = "365.44"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "7.43"
$ws.Range("E21").Value = "  +2.67%  "
$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "4.93"
$ws.Range("E23").Value = "  +3.00%  "
$ws.Range("D24").Value = "2.10"
$ws.Range("E24").Value = "  +4.53%  "
$ws.Range("D25").Value = "72.67"
$ws.Range("E25").Value = "  +7.57%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "10.09"
$ws.Range("E27").Value = "  +3.19%  "
$ws.Range("E28").Value = "  +8.70%  "
$ws.Range("D29").Value = "2.776.86"
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("D30").Value = "586.71"
$ws.Range("E30").Value = "  +3.82%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").Value = "1.43"
$ws.Range("E32").Value = "  +5.68%  "
$ws.Range("D33").Value = "8.00"
$ws.Range("E33").Value = "  +6.10%  "
$ws.Range("E34").Value = "  +4.29%  "
$ws.Range("D35").Value = "0.130"
$ws.Range("E35").Value = "  +6.10%  "
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("E37").Value = "  +5.39%  "
$ws.Range("D38").Value = "160.16"
$ws.Range("E38").Value = "  +3.12%  "
$ws.Range("E39").Value = "  +7.00%  "
$ws.Range("D40").Value = "19.32"
$ws.Range("E40").Value = "  +2.67%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "5.43"
$ws.Range("E41").Value = "  +5.87%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "0.368"
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("E43").Value = "  +7.96%  "
$ws.Range("D44").Value = "17.74"
$ws.Range("E44").Value = "  +5.65%  "
$ws.Range("E45").Value = "  +15.23%  "
$ws.Range("D46").Value = "40.75"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").Value = "156.15"
$ws.Range("E48").Value = "  +2.56%  "
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("D50").Value = "22.15"
$ws.Range("E50").Value = "  +5.18%  "
$ws.Range("D51").Value = "1.72"
$ws.Range("E51").Value = "  +3.35%  "
